$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 5.5
$ws.Range("B3").Value = 6.5
$ws.Range("C3").Value = 8.5
$ws.Range("B4").Value = 0.8
$ws.Range("C5").Value = 20

$ws.Range("J7:K7").Select()
